# Admin-side "sent notification" log: three new rows of sent-notification
# activity are being recorded above the existing history, which is why the
# original row 2 (parth / 2024-04-05) slides down to row 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data row (currently row 2) down to make room for the
# three new rows being logged; formatting/styles on the row move with it.
$ws.Rows("2:4").Insert()

# Newest entry: Jignesh, logged in 18:29:2 -> out 18:30:5 on 4/8/2024
$ws.Range("A2").Value = "Jignesh"
$ws.Range("B2").Value = "jignesh1234@gmail.com"
$ws.Range("C2").Value = 45390.771060914354
$ws.Range("C2").NumberFormat = "m/d/yy"
$ws.Range("D2").Value = "18:29:2"
$ws.Range("E2").Value = "18:30:5"

# parth, logged in 18:27:0 -> out 18:28:23 on 4/8/2024
$ws.Range("A3").Value = "parth"
$ws.Range("B3").Value = "parthpatel082828@gmail.com"
$ws.Range("C3").Value = 45390.77028099537
$ws.Range("C3").NumberFormat = "m/d/yy"
$ws.Range("D3").Value = "18:27:0"
$ws.Range("E3").Value = "18:28:23"

# Jignesh again, logged in 18:26:42 -> out 18:26:49 on 4/8/2024
$ws.Range("A4").Value = "Jignesh"
$ws.Range("B4").Value = "jignesh1234@gmail.com"
$ws.Range("C4").Value = 45390.76887056713
$ws.Range("C4").NumberFormat = "m/d/yy"
$ws.Range("D4").Value = "18:26:42"
$ws.Range("E4").Value = "18:26:49"
